# Updates the "cryptos" list (prices / hourly volume %) and swaps a couple of
# rows, matching the Aug 10 2023 GitHub Actions refresh.
#
# Most "Price" values are stored as plain text in the sheet even though they
# look numeric (e.g. "1.000", "4.123", "0.00000000117"). Excel's COM layer
# auto-converts a bare numeric-looking string assigned to Range.Value into a
# real number, which would corrupt these values (losing the thousands-style
# dots / trailing zeros). To avoid that we temporarily force a Text number
# format ("@") before writing such values, then restore the default "Normal"
# style afterwards so the cell's formatting stays identical to the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )

    $looksNumeric = $Value -match '^[+-]?[0-9]*\.?[0-9]+$'

    $range = $ws.Range($Address)
    if ($looksNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $Value
        $range.Style = "Normal"
    } else {
        $range.Value = $Value
    }
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.500.91"
Set-TextValue "E2" "  -0.79%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.849.67"
Set-TextValue "E3" "  -0.36%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.9991"
Set-TextValue "E4" "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "243.10"
Set-TextValue "E5" "  -0.64%  "

# Row 6 - XRP
Set-TextValue "D6" "0.6361"
Set-TextValue "E6" "  -1.09%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.01%  "

# Row 8 - OKB
Set-TextValue "D8" "47.68"
Set-TextValue "E8" "  +0.50%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.07559"
Set-TextValue "E9" "  +0.80%  "

# Row 11 - Solana
Set-TextValue "D11" "24.22"
Set-TextValue "E11" "  -1.05%  "

# Row 12 - TRON
Set-TextValue "D12" "0.07687"
Set-TextValue "E12" "  +0.43%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.875.81"
Set-TextValue "E13" "  +1.02%  "

# Row 14 - Polkadot
Set-TextValue "D14" "5.023"
Set-TextValue "E14" "  -0.38%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.6869"
Set-TextValue "E15" "  -0.50%  "

# Row 16 - Litecoin
Set-TextValue "E16" "  -0.04%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.000009737"
Set-TextValue "E17" "  +1.08%  "

# Row 18 - WrappedliquidstakedEther2.0
Set-TextValue "D18" "2.121.21"
Set-TextValue "E18" "  +0.57%  "

# Row 19 - Uniswap
Set-TextValue "D19" "6.221"
Set-TextValue "E19" "  +2.56%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "29.537.21"
Set-TextValue "E20" "  -0.69%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "236.31"
Set-TextValue "E21" "  +0.26%  "

# Row 22 - Avalanche
Set-TextValue "D22" "12.52"
Set-TextValue "E22" "  -1.10%  "

# Row 23 - Dai
Set-TextValue "D23" "1.000"
Set-TextValue "E23" "  +0.03%  "

# Row 24 - Chainlink
Set-TextValue "D24" "7.635"
Set-TextValue "E24" "  +2.39%  "

# Row 25 - BinanceUSD
Set-TextValue "D25" "0.9999"
Set-TextValue "E25" "  -0.05%  "

# Row 26 - Monero
Set-TextValue "D26" "156.05"
Set-TextValue "E26" "  -1.64%  "

# Row 27 - Stellar
Set-TextValue "E27" "  -2.08%  "

# Row 28 - Cosmos
Set-TextValue "D28" "8.459"
Set-TextValue "E28" "  -1.02%  "

# Row 29 - EthereumClassic
Set-TextValue "E29" "  -0.91%  "

# Row 30 - PancakeSwap
Set-TextValue "D30" "1.486"
Set-TextValue "E30" "  -0.56%  "

# Row 31 - Hedera
Set-TextValue "D31" "0.05858"
Set-TextValue "E31" "  -7.38%  "

# Row 32 - Toncoin
Set-TextValue "D32" "1.273"
Set-TextValue "E32" "  -0.47%  "

# Row 33 - Filecoin
Set-TextValue "D33" "4.123"

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "4.060"
Set-TextValue "E34" "  -0.63%  "

# Row 35 - LidoDAOToken
Set-TextValue "E35" "  +0.19%  "

# Row 36 - ARBITRUM
Set-TextValue "D36" "1.171"
Set-TextValue "E36" "  -0.24%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.7169"
Set-TextValue "E37" "  -1.75%  "

# Row 38 - HuobiToken
Set-TextValue "D38" "2.594"
Set-TextValue "E38" "  -0.56%  "

# Row 39 - MXToken
Set-TextValue "D39" "2.798"
Set-TextValue "E39" "  -1.61%  "

# Row 40 - Maker
Set-TextValue "D40" "1.237.63"
Set-TextValue "E40" "  +3.14%  "

# Row 41 - VeChain
Set-TextValue "E41" "  -0.77%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.9139"
Set-TextValue "E42" "  -0.81%  "

# Row 43 - FraxShare
Set-TextValue "D43" "6.114"
Set-TextValue "E43" "  -0.55%  "

# Row 44 - was PaxDollar, now RocketPoolETH
Set-TextValue "B44" "RocketPoolETH"
Set-TextValue "C44" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D44" "2.034.06"
Set-TextValue "E44" "  +0.73%  "

# Row 45 - was RocketPoolETH, now PaxDollar
Set-TextValue "B45" "PaxDollar"
Set-TextValue "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D45" "0.9995"
Set-TextValue "E45" "  -0.07%  "

# Row 46 - Aave
Set-TextValue "D46" "67.54"
Set-TextValue "E46" "  +1.75%  "

# Row 47 - Quant
Set-TextValue "D47" "101.68"
Set-TextValue "E47" "  -0.50%  "

# Row 48 - Aptos
Set-TextValue "D48" "7.339"
Set-TextValue "E48" "  +9.68%  "

# Row 49 - TheSandbox
Set-TextValue "D49" "0.4036"

# Row 50 - EnergySwap
Set-TextValue "D50" "9.143"
Set-TextValue "E50" "  -0.80%  "

# Row 51 - was RenderToken, now BabyDogeCoin
Set-TextValue "B51" "BabyDogeCoin"
Set-TextValue "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D51" "0.00000000117"
Set-TextValue "E51" "  -2.21%  "
